# Insert a new record row at row 455 (weekly price update for
# "Feria Lagunitas de Puerto Montt - Perejil"), pushing the existing
# rows 455:465 down to 456:466.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 455:465 down by one row, inserting a fresh blank row at 455.
$ws.Rows.Item(455).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Cells.Item(455, 1).Value = 4
$ws.Cells.Item(455, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(455, 3).Value = "Los Lagos"
$ws.Cells.Item(455, 4).Value = 45239
$ws.Cells.Item(455, 5).Value = 10
$ws.Cells.Item(455, 6).Value = 100112044
$ws.Cells.Item(455, 7).Value = "Perejil"
$ws.Cells.Item(455, 8).Value = "Sin especificar"
$ws.Cells.Item(455, 9).Value = "Primera"
$ws.Cells.Item(455, 10).Value = 50
$ws.Cells.Item(455, 11).Value = 7000
$ws.Cells.Item(455, 12).Value = 7000
$ws.Cells.Item(455, 13).Value = 7000
$ws.Cells.Item(455, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(455, 15).Value = "Región Metropolitana"
$ws.Cells.Item(455, 16).Value = 2333
$ws.Cells.Item(455, 17).Value = 3
$ws.Cells.Item(455, 18).Value = "Hortaliza"
